# Add a new "target_ids" column to the target_data table/sheet, matching the
# upstream SBTi-finance-tool commit that introduced per-row target identifiers
# (T1/T2/T3) alongside the existing target_data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("target_data")

# Extend the Excel Table ("Table4") by one column so the workbook's defined
# table range grows from A1:P120 to A1:Q120 and a 17th tableColumn is created.
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null

# Header + data for the new column.
$ws.Range("Q1").Value = "target_ids"
$ws.Range("Q2").Value = "T1"
$ws.Range("Q3").Value = "T2"
$ws.Range("Q4").Value = "T3"

# Match the saved selection/view state from the authored workbook.
$ws.Activate()
$ws.Range("Q5").Select()
